$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 7, shifting existing rows 7-13 down to 9-15
$ws.Range("A7:R8").EntireRow.Insert()

# Row 7
$ws.Range("A7").Value2 = 11
$ws.Range("B7").Value2 = 'Vega Monumental Concepción'
$ws.Range("C7").Value2 = 'Bíobío'
$ws.Range("D7").Value2 = 44868
$ws.Range("E7").Value2 = 8
$ws.Range("F7").Value2 = 300000000
$ws.Range("G7").Value2 = 'Espárragos'
$ws.Range("H7").Value2 = 'Sin especificar'
$ws.Range("I7").Value2 = 'Primera'
$ws.Range("J7").Value2 = 1000
$ws.Range("K7").Value2 = 1200
$ws.Range("L7").Value2 = 1300
$ws.Range("M7").Value2 = 1250
$ws.Range("N7").Value2 = '$/kilo'
$ws.Range("O7").Value2 = 'Región del Maule'
$ws.Range("P7").Value2 = 1250
$ws.Range("Q7").Value2 = 1
$ws.Range("R7").Value2 = 'Hortaliza'

# Row 8
$ws.Range("A8").Value2 = 11
$ws.Range("B8").Value2 = 'Vega Monumental Concepción'
$ws.Range("C8").Value2 = 'Bíobío'
$ws.Range("D8").Value2 = 44868
$ws.Range("E8").Value2 = 8
$ws.Range("F8").Value2 = 300000000
$ws.Range("G8").Value2 = 'Espárragos'
$ws.Range("H8").Value2 = 'Sin especificar'
$ws.Range("I8").Value2 = 'Segunda'
$ws.Range("J8").Value2 = 200
$ws.Range("K8").Value2 = 1000
$ws.Range("L8").Value2 = 1000
$ws.Range("M8").Value2 = 1000
$ws.Range("N8").Value2 = '$/kilo'
$ws.Range("O8").Value2 = 'Región del Maule'
$ws.Range("P8").Value2 = 1000
$ws.Range("Q8").Value2 = 1
$ws.Range("R8").Value2 = 'Hortaliza'
